# Adds 2022-Q1 data:
#  - existing "总计" sheet becomes the new "2022-Q1" per-fund holdings sheet
#  - a fresh "总计" sheet is appended (a copy of the old aggregate table)
#    with a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

$xlPasteAll     = -4104
$xlPasteFormats = -4122
$xlPasteValues  = -4163

$oldTotal = $wb.Worksheets.Item("总计")

# 1) Duplicate the aggregate ("总计") sheet, placing the copy right after the
#    original. The original keeps its position/data for now (becomes
#    "2022-Q1" below); the new copy becomes the refreshed "总计".
$oldTotal.Copy([System.Reflection.Missing]::Value, $oldTotal)
$newTotal = $wb.Worksheets.Item("总计 (2)")

$oldTotal.Name = "2022-Q1"
$newTotal.Name = "总计"

# 2) On the new "总计" sheet, shift existing data rows 2-6 down to 3-7 to make
#    room for the new 2022-Q1 row at the top (values, then formats, to avoid
#    losing the column-A style during the multi-column paste). Column A is
#    a 0-based row index, so it is bumped by one rather than carried as-is.
for ($r = 6; $r -ge 2; $r--) {
    $src = $newTotal.Range("A" + $r + ":D" + $r)
    $dst = $newTotal.Range("A" + ($r + 1) + ":D" + ($r + 1))
    $src.Copy()
    $dst.PasteSpecial($xlPasteValues)
    $src.Copy()
    $dst.PasteSpecial($xlPasteFormats)
    $newTotal.Range("A" + ($r + 1)).Value = $r - 1
}

# 3) Fill in the new top row for 2022-Q1 (copy column-A styling from row 3).
$newTotal.Range("A3").Copy()
$newTotal.Range("A2").PasteSpecial($xlPasteFormats)
$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 12
$newTotal.Range("D2").Value = 6.52

# 4) Rebuild the "2022-Q1" sheet (formerly "总计") as the per-fund holdings
#    table: extend the header style to the new columns E:H, and extend the
#    column-A index style down through the new rows 7:13.
$ws = $oldTotal

$ws.Range("D1").Copy()
$ws.Range("E1:H1").PasteSpecial($xlPasteFormats)

$ws.Range("A6").Copy()
$ws.Range("A7:A13").PasteSpecial($xlPasteFormats)

# Header row
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Data rows 2-13: A (index, numeric), B (code, text), C (name, text),
# D/E/F/G (numeric-looking values kept as text), H (rank, numeric).
$rows = @(
    @(0,  "000979", "景顺长城沪港深精选股票",                         "16.46", "82.61", "8.90", "1.4649", 5),
    @(1,  "260112", "景顺长城能源基建混合",                           "16.49", "60.89", "7.93", "1.3077", 4),
    @(2,  "009098", "景顺长城价值领航两年持有期混合",                 "11.67", "75.58", "9.86", "1.1507", 3),
    @(3,  "008850", "景顺长城价值稳进三年定期开放灵活配置混合",       "17.06", "69.71", "5.58", "0.9519", 5),
    @(4,  "008715", "景顺长城价值驱动一年持有期灵活配置混合型证券投资基金", "16.83", "62.03", "5.49", "0.9240", 3),
    @(5,  "008060", "景顺长城价值边际灵活配置混合",                   "4.93",  "80.78", "9.09", "0.4481", 3),
    @(6,  "012708", "东方红中证东方红红利低波动指数A",                "6.06",  "92.33", "1.72", "0.1042", 8),
    @(7,  "008114", "天弘中证红利低波动100指数A",                    "3.16",  "92.60", "1.65", "0.0521", 8),
    @(8,  "008115", "天弘中证红利低波动100指数C",                    "2.37",  "92.60", "1.65", "0.0391", 8),
    @(9,  "012709", "东方红中证东方红红利低波动指数C",                "2.19",  "92.33", "1.72", "0.0377", 8),
    @(10, "515100", "景顺长城中证红利低波动100ETF",                   "1.25",  "97.96", "1.74", "0.0218", 8),
    @(11, "001613", "长城久祥灵活配置混合",                           "0.30",  "88.54", "4.65", "0.0140", 5)
)

$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = "'" + $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = "'" + $row[3]
    $ws.Range("E$r").Value = "'" + $row[4]
    $ws.Range("F$r").Value = "'" + $row[5]
    $ws.Range("G$r").Value = "'" + $row[6]
    $ws.Range("H$r").Value = $row[7]
    $r++
}
